# Update the "10 Inch" pizza item names to "9 Inch" on Sheet1 (rows 11-14).
# Prices and image references are unchanged; only the Item text changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A11").Value = "Pizza Margarita 9' Inch"
$ws.Range("A12").Value = "Pizza Corn 9' Inch"
$ws.Range("A13").Value = "Pizza Onion and Capsicum 9' Inch"
$ws.Range("A14").Value = "Pizza Paneer, Veggie ( Onion and Capsicum and corn) 9' Inch"

# Restore the view/selection state recorded for the sheet.
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
